$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("A4").Value = "2026 Yamaha MT-03"
$ws.Range("A6").Value = "2026 Yamaha R3"
